{"js": "const replacements = [\n  [\"57-11=46\", \"76-31=45\"],\n  [\"97-34=63\", \"46+50=96\"],\n  [\"15+10=25\", \"59+17=76\"],\n  [\"38+2=40\", \"53+2=55\"],\n  [\"33-12=21\", \"54-14=40\"],\n  [\"40-19=21\", \"9-8=1\"],\n  [\"15+37=52\", \"32+29=61\"],\n  [\"40+44=84\", \"54-8=46\"],\n  [\"1+74=75\", \"43-14=29\"],\n  [\"69-34=35\", \"8+47=55\"],\n  [\"74-29=45\", \"27+31=58\"],\n  [\"21+22=43\", \"31+65=96\"],\n  [\"24+29=53\", \"0+71=71\"],\n  [\"35+18=53\", \"72-3=69\"],\n  [\"37-3=34\", \"81+11=92\"],\n  [\"75-21=54\", \"32+45=77\"],\n  [\"88-48=40\", \"46+38=84\"],\n  [\"99-34=65\", \"40+3=43\"],\n  [\"17+68=85\", \"59-18=41\"],\n  [\"91+8=99\", \"18+20=38\"],\n  [\"16+51=67\", \"65-24=41\"],\n  [\"63-47=16\", \"44+19=63\"],\n  [\"67-52=15\", \"56-3=53\"],\n  [\"29+61=90\", \"16+55=71\"],\n  [\"57-4=53\", \"5+4=9\"],\n  [\"65+31=96\", \"81-17=64\"],\n  [\"16-7=9\", \"61+9=70\"],\n  [\"64+4=68\", \"11+74=85\"],\n  [\"35+15=50\", \"94-53=41\"],\n  [\"80+9=89\", \"55-48=7\"],\n  [\"35+45=80\", \"4+60=64\"],\n  [\"51+36=87\", \"72+16=88\"],\n  [\"78-53=25\", \"61+25=86\"],\n  [\"73-49=24\", \"78-6=72\"],\n  [\"87-10=77\", \"80-45=35\"],\n  [\"51+21=72\", \"4+26=30\"],\n  [\"51+11=62\", \"30+20=50\"],\n  [\"26+42=68\", \"13+67=80\"],\n  [\"90-55=35\", \"74-65=9\"],\n  [\"81-53=28\", \"12+47=59\"],\n  [\"21+60=81\", \"20+74=94\"],\n  [\"34+11=45\", \"62-53=9\"],\n  [\"94-3=91\", \"31+25=56\"],\n  [\"18+13=31\", \"22+39=61\"],\n  [\"67-41=26\", \"88-44=44\"],\n  [\"36-27=9\", \"18+65=83\"],\n  [\"84-59=25\", \"6+56=62\"],\n  [\"16+74=90\", \"36-11=25\"],\n  [\"24+42=66\", \"17-10=7\"],\n  [\"44-35=9\", \"69+19=88\"],\n  [\"21-0=21\", \"73+8=81\"],\n  [\"6+42=48\", \"67-43=24\"],\n  [\"68+5=73\", \"53-36=17\"],\n  [\"15+0=15\", \"82-29=53\"],\n  [\"57-27=30\", \"63-11=52\"],\n  [\"96-70=26\", \"92-35=57\"],\n  [\"90-52=38\", \"18+49=67\"],\n  [\"76+5=81\", \"6+71=77\"],\n  [\"15+62=77\", \"79-64=15\"],\n  [\"38+26=64\", \"93+1=94\"],\n  [\"66+5=71\", \"46+42=88\"],\n  [\"78-16=62\", \"83+2=85\"],\n  [\"15+64=79\", \"93-14=79\"],\n  [\"75-7=68\", \"78-37=41\"],\n  [\"75+14=89\", \"12+28=40\"],\n  [\"20+13=33\", \"89-3=86\"],\n  [\"67+27=94\", \"57+16=73\"],\n  [\"40-2=38\", \"12-4=8\"],\n  [\"88+0=88\", \"22+53=75\"],\n  [\"23+68=91\", \"55-2=53\"],\n  [\"83-70=13\", \"88-11=77\"],\n  [\"11+61=72\", \"55+15=70\"],\n  [\"42+23=65\", \"53-2=51\"],\n  [\"65-62=3\", \"58-19=39\"],\n  [\"60-26=34\", \"71+12=83\"],\n  [\"10+82=92\", \"35+12=47\"],\n  [\"32+41=73\", \"5+85=90\"],\n  [\"68+2=70\", \"66-42=24\"],\n  [\"76+14=90\", \"28-11=17\"],\n  [\"24+48=72\", \"5+37=42\"],\n  [\"75-19=56\", \"5+28=33\"],\n  [\"62-11=51\", \"39+56=95\"],\n  [\"70-30=40\", \"6+53=59\"],\n  [\"26+18=44\", \"2+57=59\"],\n  [\"23+37=60\", \"94-13=81\"],\n  [\"20+19=39\", \"29+9=38\"],\n  [\"77-67=10\", \"5+52=57\"],\n  [\"96-93=3\", \"96-0=96\"],\n  [\"27+34=61\", \"40+41=81\"],\n  [\"93-38=55\", \"68+8=76\"],\n  [\"64+17=81\", \"82-69=13\"],\n  [\"78+5=83\", \"78-50=28\"],\n  [\"73-35=38\", \"79+8=87\"],\n  [\"90-23=67\", \"35+28=63\"],\n  [\"32-16=16\", \"17-14=3\"],\n  [\"81-77=4\", \"71-32=39\"],\n  [\"82-52=30\", \"73-73=0\"],\n  [\"92-23=69\", \"19+1=20\"],\n  [\"15+83=98\", \"12+32=44\"],\n  [\"11-6=5\", \"99-95=4\"],\n];\n\nconst body = context.document.body;\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + before);\n  }\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nreturn 'done: ' + replacements.length + ' replacements';", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Equation([string]$beforeText, [string]$afterText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $beforeText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $afterText\n    $result = $find.Execute(\n        $find.Text,             # FindText\n        $true,                  # MatchCase\n        $true,                  # MatchWholeWord\n        $false,                 # MatchWildcards\n        $false,                 # MatchSoundsLike\n        $false,                 # MatchAllWordForms\n        $true,                  # Forward\n        1,                      # Wrap (wdFindContinue)\n        $false,                 # Format\n        $find.Replacement.Text, # ReplaceWith\n        2                       # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"Replacement failed for: $beforeText\"\n    }\n}\n\nReplace-Equation \"57-11=46\" \"76-31=45\"\nReplace-Equation \"97-34=63\" \"46+50=96\"\nReplace-Equation \"15+10=25\" \"59+17=76\"\nReplace-Equation \"38+2=40\" \"53+2=55\"\nReplace-Equation \"33-12=21\" \"54-14=40\"\nReplace-Equation \"40-19=21\" \"9-8=1\"\nReplace-Equation \"15+37=52\" \"32+29=61\"\nReplace-Equation \"40+44=84\" \"54-8=46\"\nReplace-Equation \"1+74=75\" \"43-14=29\"\nReplace-Equation \"69-34=35\" \"8+47=55\"\nReplace-Equation \"74-29=45\" \"27+31=58\"\nReplace-Equation \"21+22=43\" \"31+65=96\"\nReplace-Equation \"24+29=53\" \"0+71=71\"\nReplace-Equation \"35+18=53\" \"72-3=69\"\nReplace-Equation \"37-3=34\" \"81+11=92\"\nReplace-Equation \"75-21=54\" \"32+45=77\"\nReplace-Equation \"88-48=40\" \"46+38=84\"\nReplace-Equation \"99-34=65\" \"40+3=43\"\nReplace-Equation \"17+68=85\" \"59-18=41\"\nReplace-Equation \"91+8=99\" \"18+20=38\"\nReplace-Equation \"16+51=67\" \"65-24=41\"\nReplace-Equation \"63-47=16\" \"44+19=63\"\nReplace-Equation \"67-52=15\" \"56-3=53\"\nReplace-Equation \"29+61=90\" \"16+55=71\"\nReplace-Equation \"57-4=53\" \"5+4=9\"\nReplace-Equation \"65+31=96\" \"81-17=64\"\nReplace-Equation \"16-7=9\" \"61+9=70\"\nReplace-Equation \"64+4=68\" \"11+74=85\"\nReplace-Equation \"35+15=50\" \"94-53=41\"\nReplace-Equation \"80+9=89\" \"55-48=7\"\nReplace-Equation \"35+45=80\" \"4+60=64\"\nReplace-Equation \"51+36=87\" \"72+16=88\"\nReplace-Equation \"78-53=25\" \"61+25=86\"\nReplace-Equation \"73-49=24\" \"78-6=72\"\nReplace-Equation \"87-10=77\" \"80-45=35\"\nReplace-Equation \"51+21=72\" \"4+26=30\"\nReplace-Equation \"51+11=62\" \"30+20=50\"\nReplace-Equation \"26+42=68\" \"13+67=80\"\nReplace-Equation \"90-55=35\" \"74-65=9\"\nReplace-Equation \"81-53=28\" \"12+47=59\"\nReplace-Equation \"21+60=81\" \"20+74=94\"\nReplace-Equation \"34+11=45\" \"62-53=9\"\nReplace-Equation \"94-3=91\" \"31+25=56\"\nReplace-Equation \"18+13=31\" \"22+39=61\"\nReplace-Equation \"67-41=26\" \"88-44=44\"\nReplace-Equation \"36-27=9\" \"18+65=83\"\nReplace-Equation \"84-59=25\" \"6+56=62\"\nReplace-Equation \"16+74=90\" \"36-11=25\"\nReplace-Equation \"24+42=66\" \"17-10=7\"\nReplace-Equation \"44-35=9\" \"69+19=88\"\nReplace-Equation \"21-0=21\" \"73+8=81\"\nReplace-Equation \"6+42=48\" \"67-43=24\"\nReplace-Equation \"68+5=73\" \"53-36=17\"\nReplace-Equation \"15+0=15\" \"82-29=53\"\nReplace-Equation \"57-27=30\" \"63-11=52\"\nReplace-Equation \"96-70=26\" \"92-35=57\"\nReplace-Equation \"90-52=38\" \"18+49=67\"\nReplace-Equation \"76+5=81\" \"6+71=77\"\nReplace-Equation \"15+62=77\" \"79-64=15\"\nReplace-Equation \"38+26=64\" \"93+1=94\"\nReplace-Equation \"66+5=71\" \"46+42=88\"\nReplace-Equation \"78-16=62\" \"83+2=85\"\nReplace-Equation \"15+64=79\" \"93-14=79\"\nReplace-Equation \"75-7=68\" \"78-37=41\"\nReplace-Equation \"75+14=89\" \"12+28=40\"\nReplace-Equation \"20+13=33\" \"89-3=86\"\nReplace-Equation \"67+27=94\" \"57+16=73\"\nReplace-Equation \"40-2=38\" \"12-4=8\"\nReplace-Equation \"88+0=88\" \"22+53=75\"\nReplace-Equation \"23+68=91\" \"55-2=53\"\nReplace-Equation \"83-70=13\" \"88-11=77\"\nReplace-Equation \"11+61=72\" \"55+15=70\"\nReplace-Equation \"42+23=65\" \"53-2=51\"\nReplace-Equation \"65-62=3\" \"58-19=39\"\nReplace-Equation \"60-26=34\" \"71+12=83\"\nReplace-Equation \"10+82=92\" \"35+12=47\"\nReplace-Equation \"32+41=73\" \"5+85=90\"\nReplace-Equation \"68+2=70\" \"66-42=24\"\nReplace-Equation \"76+14=90\" \"28-11=17\"\nReplace-Equation \"24+48=72\" \"5+37=42\"\nReplace-Equation \"75-19=56\" \"5+28=33\"\nReplace-Equation \"62-11=51\" \"39+56=95\"\nReplace-Equation \"70-30=40\" \"6+53=59\"\nReplace-Equation \"26+18=44\" \"2+57=59\"\nReplace-Equation \"23+37=60\" \"94-13=81\"\nReplace-Equation \"20+19=39\" \"29+9=38\"\nReplace-Equation \"77-67=10\" \"5+52=57\"\nReplace-Equation \"96-93=3\" \"96-0=96\"\nReplace-Equation \"27+34=61\" \"40+41=81\"\nReplace-Equation \"93-38=55\" \"68+8=76\"\nReplace-Equation \"64+17=81\" \"82-69=13\"\nReplace-Equation \"78+5=83\" \"78-50=28\"\nReplace-Equation \"73-35=38\" \"79+8=87\"\nReplace-Equation \"90-23=67\" \"35+28=63\"\nReplace-Equation \"32-16=16\" \"17-14=3\"\nReplace-Equation \"81-77=4\" \"71-32=39\"\nReplace-Equation \"82-52=30\" \"73-73=0\"\nReplace-Equation \"92-23=69\" \"19+1=20\"\nReplace-Equation \"15+83=98\" \"12+32=44\"\nReplace-Equation \"11-6=5\" \"99-95=4\"\n\nWrite-Output \"done\"\n"}
